$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-CellText "D2" '60.777.76'
Set-CellText "E2" '  -0.26%  '
Set-CellText "D3" '2.909.48'
Set-CellText "E3" '  -0.45%  '
Set-CellText "D4" '1.00'
Set-CellText "E4" '  -0.01%  '
Set-CellText "D5" '588.89'
Set-CellText "E5" '  +0.37%  '
Set-CellText "D6" '144.43'
Set-CellText "E6" '  -0.73%  '
Set-CellText "E7" '  +0.00%  '
Set-CellText "E8" '  -0.22%  '
Set-CellText "D9" '6.89'
Set-CellText "E9" '  +0.81%  '
Set-CellText "E11" '  -2.36%  '
Set-CellText "E12" '  -0.86%  '
Set-CellText "D13" '33.36'
Set-CellText "E13" '  -0.87%  '
Set-CellText "E14" '  -0.18%  '
Set-CellText "D15" '3.388.91'
Set-CellText "E15" '  -0.43%  '
Set-CellText "D16" '60.679.77'
Set-CellText "E16" '  -0.34%  '
Set-CellText "D17" '6.67'
Set-CellText "E17" '  -1.59%  '
Set-CellText "D18" '2.907.15'
Set-CellText "E18" '  -0.52%  '
Set-CellText "D19" '431.12'
Set-CellText "E19" '  +0.41%  '
Set-CellText "D20" '13.33'
Set-CellText "E20" '  -2.17%  '
Set-CellText "D21" '0.676'
Set-CellText "E21" '  -1.15%  '
Set-CellText "D22" '7.08'
Set-CellText "E22" '  -0.77%  '
Set-CellText "D23" '81.17'
Set-CellText "E23" '  +0.55%  '
Set-CellText "D24" '10.85'
Set-CellText "E24" '  +0.93%  '
Set-CellText "E25" '  -2.90%  '
Set-CellText "D26" '11.75'
Set-CellText "E26" '  -1.96%  '
Set-CellText "E27" '  +0.01%  '
Set-CellText "D28" '2.26'
Set-CellText "E28" '  +4.27%  '
Set-CellText "E29" '  -1.02%  '
Set-CellText "D30" '6.96'
Set-CellText "E30" '  -3.50%  '
Set-CellText "D31" '26.49'
Set-CellText "E31" '  -0.76%  '
Set-CellText "D32" '0.109'
Set-CellText "E32" '  +1.77%  '
Set-CellText "D33" '0.999'
Set-CellText "E33" '  -0.05%  '
Set-CellText "D34" '0.0₃0856'
Set-CellText "E34" '  -1.76%  '
Set-CellText "E35" '  -0.68%  '
Set-CellText "D36" '5.61'
Set-CellText "E36" '  -0.96%  '
Set-CellText "D37" '2.97'
Set-CellText "E37" '  -1.68%  '
Set-CellText "E38" '  -1.68%  '
Set-CellText "E39" '  -4.79%  '
Set-CellText "D40" '8.53'
Set-CellText "E40" '  -1.56%  '
Set-CellText "D41" '41.29'
Set-CellText "E41" '  +0.90%  '
Set-CellText "E42" '  -5.71%  '
Set-CellText "D43" '374.91'
Set-CellText "E43" '  -1.26%  '
Set-CellText "D44" '2.693.05'
Set-CellText "D45" '0.0343'
Set-CellText "E45" '  -2.84%  '
Set-CellText "D46" '133.84'
Set-CellText "E46" '  +0.96%  '
Set-CellText "D48" '23.68'
Set-CellText "E48" '  -3.40%  '
Set-CellText "E49" '  -1.01%  '
Set-CellText "E50" '  -3.44%  '
Set-CellText "E51" '  -1.15%  '
